$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.316.23'
$ws.Range("E2").Value = '  +2.40%  '
$ws.Range("D3").Value = '2.421.77'
$ws.Range("E3").Value = '  +3.11%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''556.02'
$ws.Range("E5").Value = '  +2.10%  '
$ws.Range("D6").Value = '''143.20'
$ws.Range("E6").Value = '  +4.59%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.534'
$ws.Range("E8").Value = '  +1.73%  '
$ws.Range("D9").Value = '2.421.58'
$ws.Range("E9").Value = '  +3.17%  '
$ws.Range("D10").Value = '''0.110'
$ws.Range("E10").Value = '  +4.25%  '
$ws.Range("E11").Value = '  -0.58%  '
$ws.Range("E12").Value = '  +1.59%  '
$ws.Range("D13").Value = '''0.352'
$ws.Range("E13").Value = '  +2.54%  '
$ws.Range("D14").Value = '''26.25'
$ws.Range("E14").Value = '  +6.34%  '
$ws.Range("D15").Value = '''0.0000174'
$ws.Range("E15").Value = '  +8.44%  '
$ws.Range("D16").Value = '2.859.93'
$ws.Range("E16").Value = '  +3.19%  '
$ws.Range("D17").Value = '62.202.33'
$ws.Range("E17").Value = '  +1.88%  '
$ws.Range("D18").Value = '2.420.99'
$ws.Range("E18").Value = '  +3.11%  '
$ws.Range("D19").Value = '''11.09'
$ws.Range("E19").Value = '  +4.27%  '
$ws.Range("D20").Value = '''4.21'
$ws.Range("E20").Value = '  +2.29%  '
$ws.Range("D21").Value = '''324.46'
$ws.Range("E21").Value = '  +1.77%  '
$ws.Range("D22").Value = '''6.72'
$ws.Range("E22").Value = '  +2.50%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("E24").Value = '  +5.08%  '
$ws.Range("D25").Value = '''64.94'
$ws.Range("E25").Value = '  +2.45%  '
$ws.Range("D26").Value = '''9.13'
$ws.Range("E26").Value = '  +8.63%  '
$ws.Range("D27").Value = '''575.74'
$ws.Range("E27").Value = '  +15.88%  '
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D31").Value = '0.0₃0935'
$ws.Range("E31").Value = '  +8.83%  '
$ws.Range("E32").Value = '  +6.01%  '
$ws.Range("E33").Value = '  +1.91%  '
$ws.Range("E34").Value = '  +4.15%  '
$ws.Range("E35").Value = '  +3.86%  '
$ws.Range("E36").Value = '  +8.58%  '
$ws.Range("D37").Value = '''0.999'
$ws.Range("D38").Value = '''4.83'
$ws.Range("E38").Value = '  +4.93%  '
$ws.Range("D39").Value = '''0.385'
$ws.Range("E39").Value = '  +2.23%  '
$ws.Range("E40").Value = '  +4.00%  '
$ws.Range("D41").Value = '''18.75'
$ws.Range("E41").Value = '  +1.38%  '
$ws.Range("D42").Value = '''149.31'
$ws.Range("E42").Value = '  +4.31%  '
$ws.Range("D44").Value = '''41.71'
$ws.Range("E44").Value = '  +2.70%  '
$ws.Range("E45").Value = '  +13.69%  '
$ws.Range("D46").Value = '''150.97'
$ws.Range("E46").Value = '  +5.78%  '
$ws.Range("D47").Value = '''3.64'
$ws.Range("E47").Value = '  +2.12%  '
$ws.Range("D48").Value = '''0.0541'
$ws.Range("E48").Value = '  +5.27%  '
$ws.Range("D49").Value = '''20.38'
$ws.Range("E49").Value = '  +6.71%  '
$ws.Range("E50").Value = '  +3.80%  '
$ws.Range("D51").Value = '''0.0917'
$ws.Range("E51").Value = '  +1.71%  '
